$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = 45980
$ws.Range("B31").Value = 604
$ws.Range("C31").Value = 20
$ws.Range("D31").Value = 584

$ws.Range("A31:D31").Select()
